$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update a few aggregate metrics now that trade #5 closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = 0.4   # Total P&L %
$summary.Range("B6").Value = 5     # Total Trades
$summary.Range("B9").Value = 40    # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": MarketMaking row now has 5 trades / 40% win rate.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D6").Value = 5      # Trades
$status.Range("G6").Value = 40     # Win Rate %

# ---------------------------------------------------------------------------
# New closed trade (#5) appended to both the "All Trades" and the
# "MarketMaking" strategy log sheets - identical row in each.
# ---------------------------------------------------------------------------
function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value = 5                 # Trade #

    # Date / Time are textual (e.g. "2026-02-17") and must not be
    # auto-coerced into Excel date/time serials, so force text format,
    # assign, then restore the Normal style (keeps the default xf / no
    # stray numFmt on the cell).
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"       # Date
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "23:52:18"         # Time
    $ws.Cells.Item($row, 3).Style = "Normal"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"     # Strategy
    $ws.Cells.Item($row, 5).Value = "UP"               # Side
    $ws.Cells.Item($row, 6).Value = 0.91               # Entry Price
    $ws.Cells.Item($row, 7).Value = 0.91               # Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"           # Status
    $ws.Cells.Item($row, 9).Value = 0                  # P&L %
    $ws.Cells.Item($row, 10).Value = 0                 # P&L $
    $ws.Cells.Item($row, 11).Value = 100.1             # Capital After
    $ws.Cells.Item($row, 12).Value = 0                 # Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                 # Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6               # Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"      # Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.14              # Duration (min)
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 6

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 6
